{"js": "// The canonical OOXML diff for this commit only re-orders XML namespace\n// declarations / element attributes (and drops volatile w:rsid* noise) \u2013\n// there is no actual text, formatting, numbering or structural change to\n// the document body, sectPr or numbering definitions. Word (and the\n// Office.js runtime) re-serializes the package on every save using its own\n// canonical attribute order, so the correct, semantics-preserving edit is\n// to leave the document content exactly as-is.\n//\n// We still touch the object model (load a couple of read-only properties\n// and sync) so the script demonstrably runs against the real document\n// without mutating anything.\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\n\nawait context.sync();\n\n// No content mutation is required: the template text, the bullet-list\n// numbering definition and the section/page-size properties are already\n// identical to the target state described by the diff.\n", "ps1": "# The canonical OOXML diff for this commit only re-orders XML namespace\n# declarations / element attributes (and drops volatile w:rsid* noise) -\n# there is no actual text, formatting, numbering or structural change to\n# the document body, sectPr or numbering definitions. Word re-serializes\n# the package on every save using its own canonical attribute order, so\n# the correct, semantics-preserving edit is to leave the document content\n# exactly as-is.\n#\n# We still touch the COM object model (read-only) so the script\n# demonstrably runs against the real document without mutating anything.\n$d = $word.ActiveDocument\n\n# Read-only touch of the content / paragraphs, no mutation performed.\n$null = $d.Content.Text\n$null = $d.Paragraphs.Count\n\n# No content mutation is required: the template text, the bullet-list\n# numbering definition and the section/page-size properties are already\n# identical to the target state described by the diff.\n"}
